# Swap the order of "dnasr281@gmail.com" and its paired recorder name
# in column G ("Recorded By") for every row where the cell value is an
# exact two-part, comma-separated list starting with "dnasr281@gmail.com".
# e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
#      "dnasr281@gmail.com, admin@admin.com" -> "admin@admin.com, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.SpecialCells(11).Row  # xlCellTypeLastCell = 11
if ($lastRow -lt 157) { $lastRow = 157 }

$target = "dnasr281@gmail.com"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($null -ne $val -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ", "
        if ($parts.Length -eq 2 -and $parts[0] -eq $target) {
            $cell.Value2 = "$($parts[1]), $($parts[0])"
        }
    }
}
